$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column D
$updates = @{
    3  = -7.019899999999991
    14 = -7.764500000000001
    21 = -8.776499999999997
    23 = -7.250999999999994
    25 = -8.512699999999997
    26 = -8.743199999999996
    29 = -7.044699999999997
    53 = -6.482699999999997
    57 = -8.6035
    59 = -8.173899999999998
    69 = -7.185199999999996
    79 = -6.188299999999999
    83 = -9.041399999999998
    91 = -6.579399999999998
    93 = -6.738299999999996
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
